$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the G6 input value (drives F6/H6/I6 recalculation)
# ---------------------------------------------------------------------------
$ws.Range("G6").Value = 2717

# ---------------------------------------------------------------------------
# 2. Extend the shared "F" formula pattern down into rows 7 and 8
# ---------------------------------------------------------------------------
$ws.Range("F7").Formula = "=G7/9*8"
$ws.Range("F8").Formula = "=G8/9*8"

# ---------------------------------------------------------------------------
# 3. Populate the new G7/G8 cells that mirror the prior row's H/I results
# ---------------------------------------------------------------------------
$ws.Range("G7").Formula = "=H6"
$ws.Range("G8").Formula = "=I6"

# ---------------------------------------------------------------------------
# 4. Apply the "Calculation" cell style to H5:I6 and the new G7:G8 cells
# ---------------------------------------------------------------------------
$ws.Range("H5").Style = "Calculation"
$ws.Range("I5").Style = "Calculation"
$ws.Range("H6").Style = "Calculation"
$ws.Range("I6").Style = "Calculation"
$ws.Range("G7").Style = "Calculation"
$ws.Range("G8").Style = "Calculation"

# ---------------------------------------------------------------------------
# 5. Grow the K/L "yearbook" table downward from row 255 to row 288,
#    continuing the existing L = previous + baseHeight / K = L/9*8 pattern
# ---------------------------------------------------------------------------
for ($r = 256; $r -le 288; $r++) {
    $prev = $r - 1
    $ws.Range("L$r").Formula = "=L$prev+baseHeight"
    $ws.Range("K$r").Formula = "=L$r/9*8"
}
$ws.Range("K255:L255").Copy()
$ws.Range("K256:L288").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6. Update the saved selection/active cell to reflect the new working area
# ---------------------------------------------------------------------------
$ws.Range("G7").Select()
